$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mirror the formatting of the previous data row (A10:AQ10) onto the new
# row 11 so the new year-label cell (A11) picks up the same bold/centered/
# bordered style used by every other year cell in column A, while the
# data cells keep the default (unstyled) look.
$ws.Range("A10:AQ10").Copy()
$ws.Range("A11:AQ11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "2021年"

# E11 intentionally has no value (stays an empty cell), same as the source.

$values = [ordered]@{
    "B11"  = 147.74
    "C11"  = 27.45
    "D11"  = 8.359999999999999
    "F11"  = 209.28
    "G11"  = 259.92
    "H11"  = 45.63
    "I11"  = 104.99
    "J11"  = 42.28
    "K11"  = 41.66
    "L11"  = 33.12
    "M11"  = 1.89
    "N11"  = 65.5
    "O11"  = 169.18
    "P11"  = 14.91
    "Q11"  = 50.67
    "R11"  = 146.58
    "S11"  = 14.33
    "T11"  = 155.19
    "U11"  = 0.11
    "V11"  = 72.81999999999999
    "W11"  = 23.45
    "X11"  = 234.73
    "Y11"  = 253.7
    "Z11"  = 49.8
    "AA11" = 110.12
    "AB11" = 0.33
    "AC11" = 3913.12
    "AD11" = 188.84
    "AE11" = 65.13
    "AF11" = 136.34
    "AG11" = 185.4
    "AH11" = 59
    "AI11" = 51.34
    "AJ11" = 4.34
    "AK11" = 230.83
    "AL11" = 29.52
    "AM11" = 365.49
    "AN11" = 27.67
    "AO11" = 69.93000000000001
    "AP11" = 173.69
    "AQ11" = 41.85
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
